# Change chain task data type on the "task" sheet.
# C column holds chain_task_id values. Previously some rows stored the
# "none" shared string in this column; switch those to numeric values
# (matching the numeric type already used by rows 3 and 4), and update
# row 5's other task fields to their new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("task")
$ws.Activate()

# Row 2: chain_task_id becomes numeric -1 (was text "none")
$ws.Range("C2").Value = -1

# Row 5: id, time, chain_task_id become numeric; chain_type becomes "finish"
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = -1
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = "finish"

# Update the active selection to match the saved view state
$ws.Range("C10").Select()
